# FooterPageTestData.xlsx — "Committing files after merging master"
#
# Reset the stale selection/active-tab state left over on the other sheets,
# add a new "expected Url" header label on the first ("footer") sheet, and
# leave that sheet selected at D1 as the active tab.

$wb = $excel.ActiveWorkbook

# footerSecondColumn / footerLastColumn / footerHeading: clear out their old
# lingering cell selections (H10, A2:B2, J7) back to the default A1 so they
# are no longer the active/tab-selected sheet.
$wsSecond = $wb.Worksheets.Item(2)
$wsSecond.Range("A1").Select()

$wsLast = $wb.Worksheets.Item(3)
$wsLast.Range("A1").Select()

$wsHeading = $wb.Worksheets.Item(4)
$wsHeading.Range("A1").Select()

# footer: relabel D1 from "expectedUrl" to "expected Url" and make this the
# active sheet/cell.
$wsFooter = $wb.Worksheets.Item(1)
$wsFooter.Range("D1").Value = "expected Url"
$wsFooter.Range("D1").Select()
